# New PO forecast model
# Updates the three data sheets to reflect the refreshed forecast run:
#  - "Weekly Quantity": append a new weekly bucket (row 13)
#  - "Monthly Trend": append a new monthly bucket (row 11)
#  - "PO Forecast": shift the trailing forecast window forward by 4 weeks
#                   and append one more forecast point at the end (row 21)

$wb = $excel.ActiveWorkbook

# --- Sheet: Weekly Quantity -------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A13").Value = 45676.99999999999
$wsWeekly.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeekly.Range("B13").Value = 1

# --- Sheet: Monthly Trend ----------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A11").Value = 45688.99999999999
$wsMonthly.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMonthly.Range("B11").Value = 1

# --- Sheet: PO Forecast -------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Shift the trailing forecast window (rows 13..20) forward by 28 days (4 weeks),
# i.e. the model re-ran and produced a later set of forecast dates.
$forecastRows = @(13, 14, 15, 16, 17, 18, 19, 20)
$forecastDates = @(45676.99999999999, 45683.99999999999, 45690.99999999999, 45697.99999999999, 45704.99999999999, 45711.99999999999, 45718.99999999999, 45725.99999999999)
for ($i = 0; $i -lt $forecastRows.Length; $i++) {
    $wsForecast.Cells.Item($forecastRows[$i], 1).Value = $forecastDates[$i]
}

# Append the new forecast point as row 21
$wsForecast.Range("A21").Value = 45732.99999999999
$wsForecast.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsForecast.Range("B21").Value = 2
